# append manually added monographs
# Adds three new rows of data (qid + status columns, plus QID suffixes on
# editor/author names) for monographs that were added manually with Vanderbot:
#   row 4  -> "In the Shadow of Bezalel"
#   row 10 -> "Loss and Hope: Global, Interreligious and Interdisciplinary Perspectives"
#   row 12 -> "Conundrums in Practical Theology" (also gets a hyperlink on its DOI/URL cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# --- Row 4: "In the Shadow of Bezalel" ---------------------------------
# A4/B4 are brand-new cells; copy the plain "status-column" formatting
# already used by B1 so no stray styles get created, then fill in values.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = "Q111088301"
$ws.Range("B4").Value = "added manually with Vanderbot"
$ws.Range("E4").Value = $nbsp + "Alejandro F. Botta. Q55990040"

# --- Row 10: "Loss and Hope: Global, Interreligious and Interdisciplinary Perspectives" ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = "Q111088304"
$ws.Range("B10").Value = "added manually with Vanderbot"
$ws.Range("E10").Value = "Peter Admirand Q107464477"
# A10 uses a distinct 10pt Helvetica font not used anywhere else in the sheet
$ws.Range("A10").Font.Name = "Helvetica"
$ws.Range("A10").Font.Size = 10

# --- Row 12: "Conundrums in Practical Theology" ------------------------
# A12/B12 reuse the qid-column formatting already used by A2/A5/A6/A8/A9.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = "Q111088307"
$ws.Range("B12").Value = "added manually with Vanderbot"
$ws.Range("E12").Value = "Joyce Ann Mercer Q57435308, Bonnie J. Miller-McLemore Q63038665"
$ws.Hyperlinks.Add($ws.Range("G12"), "https://ebookcentral.proquest.com/lib/vand/detail.action?docID=4694024") | Out-Null

# Restore the clipboard/selection the author ended up with.
$ws.Range("E12").Select() | Out-Null
